# Update crypto price & 1h-volume columns (D, E) for rows 2..51
# to reflect the latest scrape from the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.346.40"
$ws.Range("E2").Value = "  -6.73%  "
$ws.Range("D3").Value = "2.598.62"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'300.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "'96.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("D7").Value = "'0.576"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.56%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.557"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("E10").Value = "  -6.95%  "
$ws.Range("E11").Value = "  -3.88%  "
$ws.Range("D12").Value = "'7.80"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.49%  "
$ws.Range("D13").Value = "2.990.54"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "2.589.40"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "'0.890"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("E17").Value = "  -4.54%  "
$ws.Range("D18").Value = "43.379.06"
$ws.Range("E18").Value = "  -6.96%  "
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").Value = "0.0₃0975"
$ws.Range("E20").Value = "  -4.40%  "
$ws.Range("D21").Value = "'12.28"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.54%  "
$ws.Range("D22").Value = "'72.61"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").Value = "'264.70"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.57%  "
$ws.Range("E24").Value = "  -4.00%  "
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").Value = "'29.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'10.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").Value = "'37.52"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("D31").Value = "'6.03"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.06%  "
$ws.Range("D32").Value = "'3.59"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").Value = "'2.22"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'151.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").Value = "'0.0809"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.86%  "
$ws.Range("D37").Value = "'0.117"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.18%  "
$ws.Range("D38").Value = "'24.56"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.79%  "
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").Value = "'16.58"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").Value = "'3.55"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").Value = "'0.0313"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("E43").Value = "  -5.59%  "
$ws.Range("D44").Value = "2.040.77"
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("D45").Value = "'0.998"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'87.91"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.09%  "
$ws.Range("D47").Value = "'9.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.73%  "
$ws.Range("D48").Value = "'1.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.58%  "
$ws.Range("D49").Value = "2.846.95"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'105.64"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("D51").Value = "'0.190"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.28%  "
